# Adjust the geometry (position / size, and in some cases the vertical-flip
# flag) of several connector/rectangle shapes on slide 1 of the sequence
# diagram. The new coordinates were originally produced by nudging the
# shapes' endpoints by hand in PowerPoint, so here we reproduce the exact
# same Left/Top/Width/Height (in points) that yield the target EMU values.
#
# Note: Shape.Left/Top/Width/Height are exposed as single-precision floats
# and the host truncates (floors) points*12700 down to whole EMU, so the
# point values below were deliberately chosen (carrying extra decimal
# digits) so that floor(pt * 12700) lands exactly on the target EMU value.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Set-ShapeGeometry($Name, $Left, $Top, $Width, $Height, $FlipV) {
    $shp = $s.Shapes.Item($Name)
    if ($FlipV) {
        $shp.VerticalFlip = -1
    }
    $shp.Left = $Left
    $shp.Top = $Top
    $shp.Width = $Width
    $shp.Height = $Height
}

# id=67 - Straight Arrow Connector 66
Set-ShapeGeometry "Straight Arrow Connector 66" 170.25 122.5 45.25 0.25 $true

# id=70 - Straight Arrow Connector 69
Set-ShapeGeometry "Straight Arrow Connector 69" 171.0 166.85040283203125 86.08724975585938 0.14960630238056183 $true

# id=89 - Rectangle 88
Set-ShapeGeometry "Rectangle 88" 642.9010009765625 182.00009155273438 14.187166213989258 22.478425979614258 $false

# id=94 - Straight Arrow Connector 93 (already flipped vertically)
Set-ShapeGeometry "Straight Arrow Connector 93" 170.0 182.35623168945312 473.8407287597656 0.39377954602241516 $true

# id=97 - Rectangle 96
Set-ShapeGeometry "Rectangle 96" 840.1570434570312 237.25 11.999921798706055 26.220630645751953 $false

# id=98 - Straight Arrow Connector 97
Set-ShapeGeometry "Straight Arrow Connector 97" 170.5 238.0 670.75 0.25 $true

# id=100 - Straight Arrow Connector 99
Set-ShapeGeometry "Straight Arrow Connector 99" 170.4166259765625 262.6666259765625 669.6570434570312 0.5 $false

# id=103 - Straight Arrow Connector 102
Set-ShapeGeometry "Straight Arrow Connector 102" 170.6666259765625 76.5 161.6667022705078 0.3333858549594879 $true

# id=138 - Straight Arrow Connector 137
Set-ShapeGeometry "Straight Arrow Connector 137" 170.75 317.5 286.42071533203125 0.1434645801782608 $false

# id=156 - Straight Arrow Connector 155
Set-ShapeGeometry "Straight Arrow Connector 155" 168.0 346.04010009765625 347.6667175292969 0.2932283580303192 $false

# id=165 - Straight Arrow Connector 164
Set-ShapeGeometry "Straight Arrow Connector 164" 469.8635559082031 314.2935485839844 47.18819046020508 0.0 $false

# id=185 - Rectangle 184
Set-ShapeGeometry "Rectangle 184" 590.9730224609375 450.6667175292969 11.214173316955566 22.75244140625 $false

# id=186 - Straight Arrow Connector 185 (already flipped vertically)
Set-ShapeGeometry "Straight Arrow Connector 185" 169.74102783203125 451.0 422.09228515625 0.011653543449938297 $true

# id=188 - TextBox 187
Set-ShapeGeometry "TextBox 187" 366.0 430.612548828125 214.9078826904297 19.387481689453125 $false

# id=203 - Rectangle 202
Set-ShapeGeometry "Rectangle 202" 993.3408203125 496.5 11.999921798706055 27.25 $false

# id=204 - Straight Arrow Connector 203
Set-ShapeGeometry "Straight Arrow Connector 203" 169.33331298828125 497.0 825.166748046875 0.25 $false

# id=206 - Straight Arrow Connector 205
Set-ShapeGeometry "Straight Arrow Connector 205" 169.20001220703125 523.2000122070312 824.1408081054688 0.16204725205898285 $false
